# Add a new "WELL_COLLAR" worksheet right after "WELL_GT" (and before "GCHRON"),
# populate it with the collar-header data, format the TD/DIP/AZIMUTH values and
# the START/FINISH dates, make it the active sheet/tab, and leave the selection
# on the FINISH date cell (B3).

$wb = $excel.ActiveWorkbook

# Insert the new sheet after WELL_GT -> it lands between WELL_GT and GCHRON.
$wellGt = $wb.Worksheets.Item("WELL_GT")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wellGt)
$ws.Name = "WELL_COLLAR"

# Fill in the labels/values in the same order they were first entered so the
# shared-string table comes out in the same sequence.
$ws.Range("A1").Value = "WELL_ID"
$ws.Range("A4").Value = "AREA"
$ws.Range("A5").Value = "DRILL_CO"
$ws.Range("B5").Value = "BOART LY"
$ws.Range("B4").Value = "Upper Flats"
$ws.Range("B1").Value = "HL_001"
$ws.Range("A6").Value = "TD"
$ws.Range("A7").Value = "DIP"
$ws.Range("A8").Value = "AZIMUTH"
$ws.Range("A2").Value = "START"
$ws.Range("A3").Value = "FINISH"

# TD / DIP / AZIMUTH -> numeric, formatted with 2 decimal places.
$ws.Range("B6").Value = 1100
$ws.Range("B6").NumberFormat = "0.00"
$ws.Range("B7").Value = -89
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B8").Value = 260
$ws.Range("B6").Copy()
$ws.Range("B8").PasteSpecial(-4122)

# START / FINISH -> dates (stored as serials 7/9/2020 and 7/12/2020).
$ws.Range("B2").Value = 44021
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B3").Value = 44024
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# Make the new sheet the active tab/selection, matching the saved view state.
$ws.Range("B3").Select() | Out-Null
$ws.Activate() | Out-Null
